$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "skdslkf"
$ws.Range("B2").Value = "skdsl"
$ws.Range("C2").Value = "skdfjl"
$ws.Range("D2").Value = "slkdf"
$ws.Range("E2").Value = "skeks"
